$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "Recorded By" list of recorder names/emails, comma-separated.
# Whenever "System" is the LAST entry in that list, move it to the front
# (i.e. reverse the order of the comma-separated entries).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $raw = $cell.Value2
    if ([string]::IsNullOrEmpty($raw)) { continue }

    $parts = $raw -split ','
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1 -and $parts[$parts.Length - 1] -eq 'System') {
        $reversed = @()
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value2 = [string]::Join(', ', $reversed)
    }
}
